$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.165.14"
$ws.Range("E2").Value = "  -0.99%  "

$ws.Range("D3").Value = "1.836.14"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("D5").Value = "'233.29"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.4685"
$ws.Range("E7").Value = "  -2.71%  "

$ws.Range("D8").Value = "'0.2701"
$ws.Range("E8").Value = "  -3.81%  "

$ws.Range("D9").Value = "'0.06273"
$ws.Range("E9").Value = "  -3.63%  "

$ws.Range("D10").Value = "1.833.89"
$ws.Range("E10").Value = "  -1.78%  "

$ws.Range("D11").Value = "'0.07405"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("E13").Value = "  -2.69%  "

$ws.Range("D14").Value = "'83.66"
$ws.Range("E14").Value = "  -4.13%  "

$ws.Range("D15").Value = "'0.6178"
$ws.Range("E15").Value = "  -4.51%  "

$ws.Range("D16").Value = "30.085.17"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "'228.16"
$ws.Range("E18").Value = "  -2.60%  "

$ws.Range("D19").Value = "'0.000007273"
$ws.Range("E19").Value = "  -3.51%  "

$ws.Range("E20").Value = "  -4.99%  "

$ws.Range("D21").Value = "2.086.16"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "'4.860"
$ws.Range("E23").Value = "  -5.68%  "

$ws.Range("D24").Value = "'5.833"
$ws.Range("E24").Value = "  -4.42%  "

$ws.Range("D25").Value = "'9.185"
$ws.Range("E25").Value = "  -1.60%  "

$ws.Range("D26").Value = "'165.24"
$ws.Range("E26").Value = "  -1.14%  "

$ws.Range("E27").Value = "  -3.65%  "

$ws.Range("E28").Value = "  -2.68%  "

$ws.Range("D29").Value = "'0.1028"
$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("D31").Value = "'4.077"
$ws.Range("E31").Value = "  -4.59%  "

$ws.Range("D32").Value = "'3.787"
$ws.Range("E32").Value = "  -5.62%  "

$ws.Range("D33").Value = "'0.04795"
$ws.Range("E33").Value = "  -3.70%  "

$ws.Range("D34").Value = "'1.136"
$ws.Range("E34").Value = "  -3.63%  "

$ws.Range("D35").Value = "'0.7076"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("D36").Value = "'2.703"
$ws.Range("E36").Value = "  -0.35%  "

$ws.Range("D37").Value = "'0.01858"
$ws.Range("E37").Value = "  -3.68%  "

$ws.Range("D38").Value = "'2.645"
$ws.Range("E38").Value = "  +0.56%  "

$ws.Range("D39").Value = "'0.8934"
$ws.Range("E39").Value = "  -2.47%  "

$ws.Range("D40").Value = "'1.933"
$ws.Range("E40").Value = "  -5.87%  "

$ws.Range("D41").Value = "'104.37"
$ws.Range("E41").Value = "  -1.72%  "

$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("D43").Value = "'5.527"
$ws.Range("E43").Value = "  -0.68%  "

$ws.Range("D44").Value = "'0.3999"
$ws.Range("E44").Value = "  -4.91%  "

$ws.Range("D45").Value = "'6.951"
$ws.Range("E45").Value = "  -4.04%  "

$ws.Range("D46").Value = "'0.1190"
$ws.Range("E46").Value = "  -3.33%  "

$ws.Range("D47").Value = "'59.64"
$ws.Range("E47").Value = "  -3.71%  "

$ws.Range("D48").Value = "'8.533"
$ws.Range("E48").Value = "  -3.37%  "

$ws.Range("D49").Value = "'32.61"
$ws.Range("E49").Value = "  -2.94%  "

$ws.Range("D50").Value = "'0.05504"
$ws.Range("E50").Value = "  -2.50%  "

$ws.Range("D51").Value = "'1.358"
$ws.Range("E51").Value = "  -5.94%  "

